$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Main data table rows 2-16: columns A (label), B, C, E ---
# (column D is filled in afterwards with formulas so the engine can
#  build the shared-formula groups in the right contiguous blocks)
$rows = @(
    @("a", 1, 12, 1),
    @("b", 2, 11, 2),
    @("c", 3, 10, 3),
    @("d", 4, 9, 4),
    @("e", 9, 4, 5),
    @("f", 6, 7, 6),
    @("g", 7, 6, 7),
    @("h", 8, 5, 8),
    @("i", 9, 4, 9),
    @("j", 10, 3, 10),
    @("k", 11, 2, 11),
    @("l", 12, 1, 12),
    @("k", 10, 3, 10),
    @("l", 11, 2, 11),
    @("k", 12, 1, 12)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}

# --- Column E gets a style with applyAlignment flagged (matches the
#     workbook's new cellXfs entry) ---
$ws.Range("E2:E16").WrapText = $false

# --- Column D formulas, written per contiguous block so the engine
#     derives the same shared-formula groupings as the target file ---
$ws.Range("D2").Formula = "=B2+E2*E2"
$ws.Range("D3:D13").Formula = "=B3+E3*E3"
$ws.Range("D14:D16").Formula = "=B14+E14*E14"

# --- Extra legend rows below the table ---
$ws.Range("A17").Value = "X"
$ws.Range("A18").Value = "Y"
$ws.Range("A19").Value = "Z"
$ws.Range("B20").Value = "X1"
$ws.Range("C20").Value = "Y1"
$ws.Range("D20").Value = "Z1"

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moved to B6 ---
[void]$ws.Range("B6").Select()
